$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlTop = -4160
$xlPasteFormats = -4122
$xlPasteAll = -4104

# ---------------------------------------------------------------------------
# Row 35: Detect Cross-Browser Issues ...
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 'Detect Cross-Browser Issues for JavaScript-Based Web Applications Based on Record/Replay'
$ws.Range("A35").VerticalAlignment = $xlTop

$ws.Range("B35").Value = 'Guoquan Wu and Meimei He and Hongyin Tang and Jun Wei'
$ws.Range("B35").VerticalAlignment = $xlTop

$ws.Range("C35").Value = 42645
$ws.Range("C34").Copy()
$ws.Range("C35").PasteSpecial($xlPasteFormats)

$ws.Range("D35").Value = 44118
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial($xlPasteFormats)

$ws.Range("E35").Value = 'Haben ein Tool entwickelt, welches automatisiert cross-browser-incompatabilities (XBI) identifiziert'
$ws.Range("E35").VerticalAlignment = $xlTop
$ws.Range("E35").WrapText = $true

$ws.Range("F34").Copy()
$ws.Range("F35").PasteSpecial($xlPasteAll)
$ws.Range("F35").Value = 3

$ws.Range("H35").Value = 'https://ieeexplore.ieee.org/abstract/document/7816456'
$ws.Range("H35").VerticalAlignment = $xlTop

$ws.Rows.Item(35).RowHeight = 43.2

# ---------------------------------------------------------------------------
# Row 36: Towards Observability with (RDF) TraceStream Processing
# ---------------------------------------------------------------------------
$ws.Range("A36").Value = 'Towards Observability with (RDF) TraceStream Processing'
$ws.Range("A36").VerticalAlignment = $xlTop

$ws.Range("B36").Value = 'Mario Scrocca'
$ws.Range("B36").VerticalAlignment = $xlTop

$ws.Range("C3").Copy()
$ws.Range("C36").PasteSpecial($xlPasteAll)
$ws.Range("C36").Value = 2018

$ws.Range("D36").Value = 44118
$ws.Range("D34").Copy()
$ws.Range("D36").PasteSpecial($xlPasteFormats)

$ws.Range("E36").VerticalAlignment = $xlTop
$ws.Range("E36").WrapText = $true

$ws.Range("F34").Copy()
$ws.Range("F36").PasteSpecial($xlPasteAll)
$ws.Range("F36").Value = 4

$ws.Range("H36").Value = 'https://www.politesi.polimi.it/bitstream/10589/144741/3/2018_12_Scrocca.pdf'
$ws.Range("H36").VerticalAlignment = $xlTop

# ---------------------------------------------------------------------------
# Row 37: A Testability and Observability Framework ...
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = 'A Testability and Observability Framework to Assure Traceability Requirements on System of Systems'

$ws.Range("B37").Value = 'Leticia Morales and Miguel Ángel Olivero and Francisco José Domínguez Mayo and  J.A. Garcia-Garcia and M. Mejías'
$ws.Range("B37").VerticalAlignment = $xlTop

$ws.Range("C3").Copy()
$ws.Range("C37").PasteSpecial($xlPasteAll)
$ws.Range("C37").Value = 2020

$ws.Range("D37").Value = 44118
$ws.Range("D34").Copy()
$ws.Range("D37").PasteSpecial($xlPasteFormats)

$ws.Range("E37").Value = 'Beschreibt wie Traceability erreicht werden kann'
$ws.Range("E37").VerticalAlignment = $xlTop
$ws.Range("E37").WrapText = $true

$ws.Range("F34").Copy()
$ws.Range("F37").PasteSpecial($xlPasteAll)
$ws.Range("F37").Value = 4

$ws.Range("G37").Value = 'PDF angefragt'
$ws.Range("G37").VerticalAlignment = $xlTop

$ws.Range("H37").Value = 'https://www.researchgate.net/publication/341994379_A_Testability_and_Observability_Framework_to_Assure_Traceability_Requirements_on_System_of_Systems'
$ws.Range("H37").VerticalAlignment = $xlTop

$ws.Rows.Item(37).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 38: [US Patent Application] ...
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = '[US Patent Application] "METHODS AND SYSTEMS FOR MICROSERVICES OBSERVABILITY AUTOMATION"'
$ws.Range("A38").VerticalAlignment = $xlTop

$ws.Range("B38").Value = 'Capital One Services LLC'
$ws.Range("B38").VerticalAlignment = $xlTop

$ws.Range("C38").Value = 43357
$ws.Range("C34").Copy()
$ws.Range("C38").PasteSpecial($xlPasteFormats)

$ws.Range("D38").Value = 44118
$ws.Range("D34").Copy()
$ws.Range("D38").PasteSpecial($xlPasteFormats)

$ws.Range("E38").Value = 'Beschreibt einen Ansatz, wie man die aktuell verschiedenen Technologien zusammenfasst und homogenisiert'
$ws.Range("E38").VerticalAlignment = $xlTop
$ws.Range("E38").WrapText = $true

$ws.Range("F34").Copy()
$ws.Range("F38").PasteSpecial($xlPasteAll)
$ws.Range("F38").Value = 3

$ws.Range("H38").Value = 'https://patentimages.storage.googleapis.com/8c/79/f2/fcc54da37c1b49/US20200092180A1.pdf'
$ws.Range("H38").VerticalAlignment = $xlTop

$ws.Rows.Item(38).RowHeight = 43.2

# ---------------------------------------------------------------------------
# View state: selection moves to H39 after the newly-entered rows
# ---------------------------------------------------------------------------
$ws.Range("H39").Select()
